$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,8).Value = "kitchens"
$ws.Cells.Item(2,9).Value = "target"
$ws.Cells.Item(2,11).Value = "j"
$ws.Cells.Item(2,12).Value = "stimuli/img_mgnmm.png"
$ws.Cells.Item(2,13).Value = 79.1470588235294
$ws.Cells.Item(2,14).Value = 60.38235294117647
$ws.Cells.Item(2,15).Value = 69.76470588235294
$ws.Cells.Item(2,16).Value = 34
$ws.Cells.Item(2,17).Value = 8
$ws.Cells.Item(2,18).Value = 8
$ws.Cells.Item(2,19).Value = 8
$ws.Cells.Item(2,20).Value = 8
$ws.Cells.Item(2,21).Value = 8
$ws.Cells.Item(2,22).Value = 8
$ws.Cells.Item(3,6).Value = 2
$ws.Cells.Item(3,12).Value = "stimuli/img_esb4r.png"
$ws.Cells.Item(3,13).Value = 60.73529411764706
$ws.Cells.Item(3,14).Value = 38.58823529411764
$ws.Cells.Item(3,15).Value = 49.66176470588235
$ws.Cells.Item(3,17).Value = 3
$ws.Cells.Item(3,18).Value = 3
$ws.Cells.Item(3,19).Value = 3
$ws.Cells.Item(3,20).Value = 3
$ws.Cells.Item(3,21).Value = 3
$ws.Cells.Item(3,22).Value = 3
$ws.Cells.Item(4,6).Value = 3
$ws.Cells.Item(4,12).Value = "stimuli/img_1ao2d.png"
$ws.Cells.Item(4,13).Value = 38.77777777777778
$ws.Cells.Item(4,14).Value = 18.75
$ws.Cells.Item(4,15).Value = 28.76388888888889
$ws.Cells.Item(4,16).Value = 36
$ws.Cells.Item(4,17).Value = 1
$ws.Cells.Item(4,18).Value = 1
$ws.Cells.Item(4,19).Value = 1
$ws.Cells.Item(4,20).Value = 1
$ws.Cells.Item(4,21).Value = 1
$ws.Cells.Item(4,22).Value = 1
$ws.Cells.Item(5,6).Value = 4
$ws.Cells.Item(5,12).Value = "stimuli/img_n9xll.png"
$ws.Cells.Item(5,13).Value = 77.14285714285714
$ws.Cells.Item(5,14).Value = 59.21428571428572
$ws.Cells.Item(5,15).Value = 68.17857142857143
$ws.Cells.Item(5,16).Value = 42
$ws.Cells.Item(6,6).Value = 5
$ws.Cells.Item(6,8).Value = "bedrooms"
$ws.Cells.Item(6,9).Value = "distractor"
$ws.Cells.Item(6,11).Value = "f"
$ws.Cells.Item(6,12).Value = "stimuli/img_dkqas.png"
$ws.Cells.Item(6,13).Value = 78.57894736842105
$ws.Cells.Item(6,14).Value = 57.71052631578947
$ws.Cells.Item(6,15).Value = 68.14473684210526
$ws.Cells.Item(6,16).Value = 38
$ws.Cells.Item(6,17).Value = 7
$ws.Cells.Item(6,18).Value = 7
$ws.Cells.Item(6,19).Value = 7
$ws.Cells.Item(6,20).Value = 7
$ws.Cells.Item(6,21).Value = 7
$ws.Cells.Item(6,22).Value = 7
$ws.Cells.Item(7,6).Value = 6
$ws.Cells.Item(7,8).Value = "living_rooms"
$ws.Cells.Item(7,9).Value = "distractor"
$ws.Cells.Item(7,11).Value = "f"
$ws.Cells.Item(7,12).Value = "stimuli/img_pna7l.png"
$ws.Cells.Item(7,13).Value = 85.53333333333333
$ws.Cells.Item(7,14).Value = 67.97777777777777
$ws.Cells.Item(7,15).Value = 76.75555555555556
$ws.Cells.Item(7,16).Value = 45
$ws.Cells.Item(7,17).Value = 9
$ws.Cells.Item(7,18).Value = 9
$ws.Cells.Item(7,19).Value = 9
$ws.Cells.Item(7,20).Value = 9
$ws.Cells.Item(7,21).Value = 9
$ws.Cells.Item(7,22).Value = 9
$ws.Cells.Item(8,6).Value = 7
$ws.Cells.Item(8,8).Value = "kitchens"
$ws.Cells.Item(8,9).Value = "target"
$ws.Cells.Item(8,11).Value = "j"
$ws.Cells.Item(8,12).Value = "stimuli/img_wppku.png"
$ws.Cells.Item(8,13).Value = 75.02941176470588
$ws.Cells.Item(8,14).Value = 53.05882352941177
$ws.Cells.Item(8,15).Value = 64.04411764705883
$ws.Cells.Item(8,16).Value = 34
$ws.Cells.Item(8,17).Value = 6
$ws.Cells.Item(8,18).Value = 6
$ws.Cells.Item(8,19).Value = 6
$ws.Cells.Item(8,20).Value = 6
$ws.Cells.Item(8,21).Value = 6
$ws.Cells.Item(8,22).Value = 6
$ws.Cells.Item(9,6).Value = 8
$ws.Cells.Item(9,8).Value = "kitchens"
$ws.Cells.Item(9,9).Value = "target"
$ws.Cells.Item(9,11).Value = "j"
$ws.Cells.Item(9,12).Value = "stimuli/img_7ucnr.png"
$ws.Cells.Item(9,13).Value = 70.39393939393939
$ws.Cells.Item(9,14).Value = 47.90909090909091
$ws.Cells.Item(9,15).Value = 59.15151515151515
$ws.Cells.Item(9,16).Value = 33
$ws.Cells.Item(9,17).Value = 5
$ws.Cells.Item(9,18).Value = 5
$ws.Cells.Item(9,19).Value = 5
$ws.Cells.Item(9,20).Value = 5
$ws.Cells.Item(9,21).Value = 5
$ws.Cells.Item(9,22).Value = 5
$ws.Cells.Item(10,6).Value = 9
$ws.Cells.Item(10,8).Value = "living_rooms"
$ws.Cells.Item(10,9).Value = "distractor"
$ws.Cells.Item(10,11).Value = "f"
$ws.Cells.Item(10,12).Value = "stimuli/img_amsgw.png"
$ws.Cells.Item(10,13).Value = 86.08510638297872
$ws.Cells.Item(10,14).Value = 65.95744680851064
$ws.Cells.Item(10,15).Value = 76.02127659574468
$ws.Cells.Item(10,16).Value = 47
$ws.Cells.Item(10,17).Value = 9
$ws.Cells.Item(10,18).Value = 9
$ws.Cells.Item(10,19).Value = 9
$ws.Cells.Item(10,20).Value = 8
$ws.Cells.Item(10,21).Value = 9
$ws.Cells.Item(10,22).Value = 8
$ws.Cells.Item(11,6).Value = 10
$ws.Cells.Item(11,12).Value = "stimuli/img_ac0ey.png"
$ws.Cells.Item(11,13).Value = 86.62222222222222
$ws.Cells.Item(11,14).Value = 70.02222222222223
$ws.Cells.Item(11,15).Value = 78.32222222222222
$ws.Cells.Item(11,16).Value = 45
$ws.Cells.Item(11,20).Value = 9
$ws.Cells.Item(11,22).Value = 9
$ws.Cells.Item(12,6).Value = 11
$ws.Cells.Item(12,12).Value = "stimuli/img_5m6x4.png"
$ws.Cells.Item(12,13).Value = 80.23076923076923
$ws.Cells.Item(12,14).Value = 58.41025641025641
$ws.Cells.Item(12,15).Value = 69.32051282051282
$ws.Cells.Item(12,16).Value = 39
$ws.Cells.Item(12,17).Value = 7
$ws.Cells.Item(12,18).Value = 7
$ws.Cells.Item(12,19).Value = 7
$ws.Cells.Item(12,20).Value = 7
$ws.Cells.Item(12,21).Value = 7
$ws.Cells.Item(12,22).Value = 7
$ws.Cells.Item(13,6).Value = 12
$ws.Cells.Item(13,12).Value = "stimuli/img_95hiv.png"
$ws.Cells.Item(13,13).Value = 84.04545454545455
$ws.Cells.Item(13,14).Value = 67.31818181818181
$ws.Cells.Item(13,15).Value = 75.68181818181819
$ws.Cells.Item(13,16).Value = 44
$ws.Cells.Item(13,20).Value = 8
$ws.Cells.Item(13,21).Value = 8
$ws.Cells.Item(13,22).Value = 8
$ws.Cells.Item(14,6).Value = 13
$ws.Cells.Item(14,12).Value = "stimuli/img_ikk62.png"
$ws.Cells.Item(14,13).Value = 37.48780487804878
$ws.Cells.Item(14,14).Value = 21.07317073170732
$ws.Cells.Item(14,15).Value = 29.28048780487805
$ws.Cells.Item(14,16).Value = 41
$ws.Cells.Item(14,17).Value = 1
$ws.Cells.Item(14,18).Value = 1
$ws.Cells.Item(14,19).Value = 1
$ws.Cells.Item(14,20).Value = 1
$ws.Cells.Item(14,21).Value = 1
$ws.Cells.Item(14,22).Value = 1
$ws.Cells.Item(15,6).Value = 14
$ws.Cells.Item(15,8).Value = "bedrooms"
$ws.Cells.Item(15,12).Value = "stimuli/img_u1rxv.png"
$ws.Cells.Item(15,13).Value = 75.63636363636364
$ws.Cells.Item(15,14).Value = 54.27272727272727
$ws.Cells.Item(15,15).Value = 64.95454545454545
$ws.Cells.Item(15,16).Value = 44
$ws.Cells.Item(15,17).Value = 6
$ws.Cells.Item(15,18).Value = 6
$ws.Cells.Item(15,19).Value = 6
$ws.Cells.Item(15,20).Value = 6
$ws.Cells.Item(15,21).Value = 6
$ws.Cells.Item(15,22).Value = 6
$ws.Cells.Item(16,6).Value = 15
$ws.Cells.Item(16,8).Value = "living_rooms"
$ws.Cells.Item(16,9).Value = "distractor"
$ws.Cells.Item(16,11).Value = "f"
$ws.Cells.Item(16,12).Value = "stimuli/img_f63yi.png"
$ws.Cells.Item(16,13).Value = 85.275
$ws.Cells.Item(16,14).Value = 68.475
$ws.Cells.Item(16,15).Value = 76.875
$ws.Cells.Item(16,16).Value = 40
$ws.Cells.Item(16,21).Value = 8
$ws.Cells.Item(17,6).Value = 16
$ws.Cells.Item(17,8).Value = "living_rooms"
$ws.Cells.Item(17,12).Value = "stimuli/img_6zz63.png"
$ws.Cells.Item(17,13).Value = 87.66666666666667
$ws.Cells.Item(17,14).Value = 70.6
$ws.Cells.Item(17,15).Value = 79.13333333333333
$ws.Cells.Item(17,16).Value = 45
$ws.Cells.Item(17,17).Value = 9
$ws.Cells.Item(17,18).Value = 10
$ws.Cells.Item(17,19).Value = 10
$ws.Cells.Item(17,20).Value = 9
$ws.Cells.Item(17,21).Value = 9
$ws.Cells.Item(17,22).Value = 9
$ws.Cells.Item(18,6).Value = 17
$ws.Cells.Item(18,8).Value = "living_rooms"
$ws.Cells.Item(18,9).Value = "distractor"
$ws.Cells.Item(18,11).Value = "f"
$ws.Cells.Item(18,12).Value = "stimuli/img_1zhz6.png"
$ws.Cells.Item(18,13).Value = 49.02272727272727
$ws.Cells.Item(18,14).Value = 32.77272727272727
$ws.Cells.Item(18,15).Value = 40.89772727272727
$ws.Cells.Item(18,16).Value = 44
$ws.Cells.Item(18,17).Value = 3
$ws.Cells.Item(18,18).Value = 3
$ws.Cells.Item(18,19).Value = 3
$ws.Cells.Item(18,20).Value = 3
$ws.Cells.Item(18,21).Value = 3
$ws.Cells.Item(18,22).Value = 3
$ws.Cells.Item(19,6).Value = 18
$ws.Cells.Item(19,8).Value = "kitchens"
$ws.Cells.Item(19,9).Value = "target"
$ws.Cells.Item(19,11).Value = "j"
$ws.Cells.Item(19,12).Value = "stimuli/img_mawe6.png"
$ws.Cells.Item(19,13).Value = 83.48387096774194
$ws.Cells.Item(19,14).Value = 65.54838709677419
$ws.Cells.Item(19,15).Value = 74.51612903225806
$ws.Cells.Item(19,16).Value = 31
$ws.Cells.Item(19,17).Value = 9
$ws.Cells.Item(19,18).Value = 9
$ws.Cells.Item(19,19).Value = 9
$ws.Cells.Item(19,20).Value = 9
$ws.Cells.Item(19,21).Value = 9
$ws.Cells.Item(19,22).Value = 9
$ws.Cells.Item(20,6).Value = 19
$ws.Cells.Item(20,8).Value = "bedrooms"
$ws.Cells.Item(20,9).Value = "distractor"
$ws.Cells.Item(20,11).Value = "f"
$ws.Cells.Item(20,12).Value = "stimuli/img_d9ogj.png"
$ws.Cells.Item(20,13).Value = 76.86842105263158
$ws.Cells.Item(20,14).Value = 53.5
$ws.Cells.Item(20,15).Value = 65.18421052631578
$ws.Cells.Item(20,16).Value = 38
$ws.Cells.Item(20,17).Value = 6
$ws.Cells.Item(20,18).Value = 6
$ws.Cells.Item(20,19).Value = 6
$ws.Cells.Item(20,20).Value = 6
$ws.Cells.Item(20,21).Value = 6
$ws.Cells.Item(20,22).Value = 6
$ws.Cells.Item(21,6).Value = 20
$ws.Cells.Item(21,12).Value = "stimuli/img_qmgwq.png"
$ws.Cells.Item(21,13).Value = 84.58333333333333
$ws.Cells.Item(21,14).Value = 64.44444444444444
$ws.Cells.Item(21,15).Value = 74.51388888888889
$ws.Cells.Item(21,16).Value = 36
$ws.Cells.Item(21,17).Value = 9
$ws.Cells.Item(21,18).Value = 9
$ws.Cells.Item(21,19).Value = 9
$ws.Cells.Item(21,20).Value = 9
$ws.Cells.Item(21,21).Value = 9
$ws.Cells.Item(21,22).Value = 9
$ws.Cells.Item(22,6).Value = 21
$ws.Cells.Item(22,12).Value = "stimuli/img_kwxq1.png"
$ws.Cells.Item(22,13).Value = 68.53125
$ws.Cells.Item(22,14).Value = 44.09375
$ws.Cells.Item(22,15).Value = 56.3125
$ws.Cells.Item(22,16).Value = 32
$ws.Cells.Item(22,17).Value = 4
$ws.Cells.Item(22,18).Value = 4
$ws.Cells.Item(22,19).Value = 4
$ws.Cells.Item(22,20).Value = 4
$ws.Cells.Item(22,21).Value = 4
$ws.Cells.Item(22,22).Value = 4
$ws.Cells.Item(23,6).Value = 22
$ws.Cells.Item(23,12).Value = "stimuli/img_8dacu.png"
$ws.Cells.Item(23,13).Value = 76.38461538461539
$ws.Cells.Item(23,14).Value = 53.64102564102564
$ws.Cells.Item(23,15).Value = 65.01282051282051
$ws.Cells.Item(23,16).Value = 39
$ws.Cells.Item(23,17).Value = 6
$ws.Cells.Item(23,18).Value = 6
$ws.Cells.Item(23,19).Value = 6
$ws.Cells.Item(23,20).Value = 6
$ws.Cells.Item(23,21).Value = 6
$ws.Cells.Item(23,22).Value = 6
$ws.Cells.Item(24,6).Value = 23
$ws.Cells.Item(24,8).Value = "kitchens"
$ws.Cells.Item(24,9).Value = "target"
$ws.Cells.Item(24,11).Value = "j"
$ws.Cells.Item(24,12).Value = "stimuli/img_zi8qc.png"
$ws.Cells.Item(24,13).Value = 77.14285714285714
$ws.Cells.Item(24,14).Value = 57.02857142857143
$ws.Cells.Item(24,15).Value = 67.08571428571429
$ws.Cells.Item(24,16).Value = 35
$ws.Cells.Item(25,6).Value = 24
$ws.Cells.Item(25,12).Value = "stimuli/img_89rmb.png"
$ws.Cells.Item(25,13).Value = 55.18518518518518
$ws.Cells.Item(25,14).Value = 29.25925925925926
$ws.Cells.Item(25,15).Value = 42.22222222222222
$ws.Cells.Item(25,16).Value = 27
$ws.Cells.Item(25,17).Value = 2
$ws.Cells.Item(25,18).Value = 2
$ws.Cells.Item(25,19).Value = 2
$ws.Cells.Item(25,20).Value = 2
$ws.Cells.Item(25,21).Value = 2
$ws.Cells.Item(25,22).Value = 2
$ws.Cells.Item(26,6).Value = 25
$ws.Cells.Item(26,12).Value = "stimuli/img_xdhz2.png"
$ws.Cells.Item(26,13).Value = 63.3
$ws.Cells.Item(26,14).Value = 37.25
$ws.Cells.Item(26,15).Value = 50.275
$ws.Cells.Item(26,16).Value = 40
$ws.Cells.Item(26,17).Value = 3
$ws.Cells.Item(26,18).Value = 3
$ws.Cells.Item(26,19).Value = 3
$ws.Cells.Item(26,20).Value = 3
$ws.Cells.Item(26,21).Value = 3
$ws.Cells.Item(26,22).Value = 3
$ws.Cells.Item(27,6).Value = 26
$ws.Cells.Item(27,12).Value = "stimuli/img_d26ik.png"
$ws.Cells.Item(27,13).Value = 77.73809523809524
$ws.Cells.Item(27,14).Value = 60.66666666666666
$ws.Cells.Item(27,15).Value = 69.20238095238095
$ws.Cells.Item(27,16).Value = 42
$ws.Cells.Item(27,17).Value = 7
$ws.Cells.Item(27,18).Value = 7
$ws.Cells.Item(27,19).Value = 7
$ws.Cells.Item(27,20).Value = 7
$ws.Cells.Item(27,21).Value = 7
$ws.Cells.Item(27,22).Value = 7
$ws.Cells.Item(28,6).Value = 27
$ws.Cells.Item(28,8).Value = "kitchens"
$ws.Cells.Item(28,9).Value = "target"
$ws.Cells.Item(28,11).Value = "j"
$ws.Cells.Item(28,12).Value = "stimuli/img_7w5tw.png"
$ws.Cells.Item(28,13).Value = 53.2258064516129
$ws.Cells.Item(28,14).Value = 28.90322580645161
$ws.Cells.Item(28,15).Value = 41.06451612903226
$ws.Cells.Item(28,16).Value = 31
$ws.Cells.Item(28,17).Value = 2
$ws.Cells.Item(28,18).Value = 2
$ws.Cells.Item(28,19).Value = 2
$ws.Cells.Item(28,20).Value = 2
$ws.Cells.Item(28,21).Value = 2
$ws.Cells.Item(28,22).Value = 2
$ws.Cells.Item(29,6).Value = 28
$ws.Cells.Item(29,12).Value = "stimuli/img_ewrjk.png"
$ws.Cells.Item(29,13).Value = 73.0909090909091
$ws.Cells.Item(29,14).Value = 53.39393939393939
$ws.Cells.Item(29,15).Value = 63.24242424242424
$ws.Cells.Item(29,16).Value = 33
$ws.Cells.Item(29,17).Value = 6
$ws.Cells.Item(29,18).Value = 6
$ws.Cells.Item(29,19).Value = 6
$ws.Cells.Item(29,20).Value = 6
$ws.Cells.Item(29,21).Value = 6
$ws.Cells.Item(29,22).Value = 6
$ws.Cells.Item(30,6).Value = 29
$ws.Cells.Item(30,12).Value = "stimuli/img_lszzj.png"
$ws.Cells.Item(30,13).Value = 64.70588235294117
$ws.Cells.Item(30,14).Value = 45.58823529411764
$ws.Cells.Item(30,15).Value = 55.14705882352941
$ws.Cells.Item(30,16).Value = 34
$ws.Cells.Item(30,17).Value = 4
$ws.Cells.Item(30,18).Value = 4
$ws.Cells.Item(30,19).Value = 4
$ws.Cells.Item(30,20).Value = 4
$ws.Cells.Item(30,21).Value = 4
$ws.Cells.Item(30,22).Value = 4
$ws.Cells.Item(31,6).Value = 30
$ws.Cells.Item(31,12).Value = "stimuli/img_mjxmq.png"
$ws.Cells.Item(31,13).Value = 77.07692307692308
$ws.Cells.Item(31,14).Value = 58.15384615384615
$ws.Cells.Item(31,15).Value = 67.61538461538461
$ws.Cells.Item(31,16).Value = 39
$ws.Cells.Item(31,17).Value = 7
$ws.Cells.Item(31,18).Value = 7
$ws.Cells.Item(31,19).Value = 7
$ws.Cells.Item(31,20).Value = 7
$ws.Cells.Item(31,21).Value = 7
$ws.Cells.Item(31,22).Value = 7
$ws.Cells.Item(32,6).Value = 31
$ws.Cells.Item(33,6).Value = 32
$ws.Cells.Item(33,12).Value = "stimuli/img_r2lxk.png"
$ws.Cells.Item(33,13).Value = 89.24242424242425
$ws.Cells.Item(33,14).Value = 67.6969696969697
$ws.Cells.Item(33,15).Value = 78.46969696969697
$ws.Cells.Item(33,16).Value = 33
$ws.Cells.Item(33,17).Value = 10
$ws.Cells.Item(33,18).Value = 10
$ws.Cells.Item(33,19).Value = 10
$ws.Cells.Item(33,20).Value = 10
$ws.Cells.Item(33,21).Value = 10
$ws.Cells.Item(33,22).Value = 10
$ws.Cells.Item(34,6).Value = 33
$ws.Cells.Item(34,12).Value = "stimuli/img_yosqb.png"
$ws.Cells.Item(34,13).Value = 50.88372093023256
$ws.Cells.Item(34,14).Value = 30.11627906976744
$ws.Cells.Item(34,15).Value = 40.5
$ws.Cells.Item(34,16).Value = 43
$ws.Cells.Item(34,17).Value = 3
$ws.Cells.Item(34,18).Value = 3
$ws.Cells.Item(34,19).Value = 3
$ws.Cells.Item(34,20).Value = 3
$ws.Cells.Item(34,21).Value = 3
$ws.Cells.Item(34,22).Value = 3
$ws.Cells.Item(35,6).Value = 34
$ws.Cells.Item(35,12).Value = "stimuli/img_0jzz7.png"
$ws.Cells.Item(35,13).Value = 84.85106382978724
$ws.Cells.Item(35,14).Value = 68.87234042553192
$ws.Cells.Item(35,15).Value = 76.86170212765958
$ws.Cells.Item(35,16).Value = 47
$ws.Cells.Item(35,20).Value = 9
$ws.Cells.Item(35,22).Value = 9
$ws.Cells.Item(36,6).Value = 35
$ws.Cells.Item(36,8).Value = "bedrooms"
$ws.Cells.Item(36,12).Value = "stimuli/img_ybbmx.png"
$ws.Cells.Item(36,13).Value = 55.24324324324324
$ws.Cells.Item(36,14).Value = 36.75675675675676
$ws.Cells.Item(36,15).Value = 46
$ws.Cells.Item(36,16).Value = 37
$ws.Cells.Item(36,17).Value = 3
$ws.Cells.Item(36,18).Value = 3
$ws.Cells.Item(36,19).Value = 3
$ws.Cells.Item(36,20).Value = 3
$ws.Cells.Item(36,21).Value = 3
$ws.Cells.Item(36,22).Value = 3
$ws.Cells.Item(37,6).Value = 36
$ws.Cells.Item(37,8).Value = "bedrooms"
$ws.Cells.Item(37,9).Value = "distractor"
$ws.Cells.Item(37,11).Value = "f"
$ws.Cells.Item(37,12).Value = "stimuli/img_fea1z.png"
$ws.Cells.Item(37,13).Value = 79.45945945945945
$ws.Cells.Item(37,14).Value = 56.24324324324324
$ws.Cells.Item(37,15).Value = 67.85135135135135
$ws.Cells.Item(37,16).Value = 37
$ws.Cells.Item(37,17).Value = 7
$ws.Cells.Item(37,18).Value = 7
$ws.Cells.Item(37,19).Value = 7
$ws.Cells.Item(37,20).Value = 7
$ws.Cells.Item(37,21).Value = 7
$ws.Cells.Item(37,22).Value = 7
$ws.Cells.Item(38,6).Value = 37
$ws.Cells.Item(38,8).Value = "bedrooms"
$ws.Cells.Item(38,12).Value = "stimuli/img_iqmdm.png"
$ws.Cells.Item(38,13).Value = 79.38888888888889
$ws.Cells.Item(38,14).Value = 58.36111111111111
$ws.Cells.Item(38,15).Value = 68.875
$ws.Cells.Item(38,16).Value = 36
$ws.Cells.Item(38,17).Value = 7
$ws.Cells.Item(38,18).Value = 7
$ws.Cells.Item(38,19).Value = 7
$ws.Cells.Item(38,20).Value = 7
$ws.Cells.Item(38,21).Value = 7
$ws.Cells.Item(38,22).Value = 7
$ws.Cells.Item(39,6).Value = 38
$ws.Cells.Item(39,8).Value = "kitchens"
$ws.Cells.Item(39,9).Value = "target"
$ws.Cells.Item(39,11).Value = "j"
$ws.Cells.Item(39,12).Value = "stimuli/img_7ed9m.png"
$ws.Cells.Item(39,13).Value = 80.71875
$ws.Cells.Item(39,14).Value = 58.65625
$ws.Cells.Item(39,15).Value = 69.6875
$ws.Cells.Item(39,16).Value = 32
$ws.Cells.Item(39,17).Value = 8
$ws.Cells.Item(39,18).Value = 8
$ws.Cells.Item(39,19).Value = 8
$ws.Cells.Item(39,20).Value = 8
$ws.Cells.Item(39,22).Value = 8
$ws.Cells.Item(40,6).Value = 39
$ws.Cells.Item(40,8).Value = "living_rooms"
$ws.Cells.Item(40,9).Value = "distractor"
$ws.Cells.Item(40,11).Value = "f"
$ws.Cells.Item(40,12).Value = "stimuli/img_wgkqa.png"
$ws.Cells.Item(40,13).Value = 87.25581395348837
$ws.Cells.Item(40,14).Value = 71.13953488372093
$ws.Cells.Item(40,15).Value = 79.19767441860465
$ws.Cells.Item(40,16).Value = 43
$ws.Cells.Item(40,17).Value = 10
$ws.Cells.Item(40,18).Value = 10
$ws.Cells.Item(40,19).Value = 10
$ws.Cells.Item(40,20).Value = 9
$ws.Cells.Item(40,21).Value = 9
$ws.Cells.Item(40,22).Value = 9
$ws.Cells.Item(41,6).Value = 40
$ws.Cells.Item(41,8).Value = "kitchens"
$ws.Cells.Item(41,9).Value = "target"
$ws.Cells.Item(41,11).Value = "j"
$ws.Cells.Item(41,12).Value = "stimuli/img_z293c.png"
$ws.Cells.Item(41,13).Value = 71.26470588235294
$ws.Cells.Item(41,14).Value = 46.88235294117647
$ws.Cells.Item(41,15).Value = 59.07352941176471
$ws.Cells.Item(41,16).Value = 34
$ws.Cells.Item(41,17).Value = 5
$ws.Cells.Item(41,18).Value = 5
$ws.Cells.Item(41,19).Value = 5
$ws.Cells.Item(41,20).Value = 5
$ws.Cells.Item(41,21).Value = 5
$ws.Cells.Item(41,22).Value = 5

Write-Host "Applied kitchens categorization row permutation update"
